$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos.xlsx stores every Price/Link/Coin/Volume cell as literal text
# (inline strings), even though many of the price strings look numeric
# (e.g. "486.99", "1.00"). A plain Range.Value assignment lets Excel
# "smart-parse" those into real numbers, which would corrupt the data,
# so each target cell is forced to the Text number format right before
# its new value is written.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.254.40'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.473.24'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.99%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '486.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.40'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +12.03%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.480.57'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +9.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0965'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.331'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.14%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.903.92'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.218.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.02'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +7.38%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.478.02'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.17%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +8.43%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '316.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.94%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +8.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.35'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.410'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +7.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.75%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +6.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.586.75'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +8.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0787'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +11.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.94'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.14'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.09%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.18'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.78%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +8.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.15'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.93%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.29%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.88%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.76%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +13.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0926'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.69%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '257.87'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +15.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.20'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.03%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0228'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.47'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.858.62'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.28%  '
